$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the full previous data region (values + formatting) before rebuilding it
$ws.Range("A1:C25").Clear()

# Column layout: column A keeps its own definition (A1:A1); column B gets its own
# width/style instead of inheriting from the old combined A:B range definition.
$ws.Columns.Item(1).ColumnWidth = 30.7109375
$ws.Columns.Item(2).ColumnWidth = 60.7109375
$ws.Columns.Item(3).ColumnWidth = 60.7109375

$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Rows.Item(1).AutoFit()

$ws.Range("B2").Value = 'LOM3106'
$ws.Range("C2").Value = 'LOM3106'
$ws.Rows.Item(2).AutoFit()

$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Ciência dos Materiais Computacional'
$ws.Range("C3").Value = ' Ciência dos Materiais Computacional'
$ws.Rows.Item(3).AutoFit()

$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Computational Materials Science'
$ws.Range("C4").Value = 'Computational Materials Science'
$ws.Rows.Item(4).AutoFit()

$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '2'
$ws.Range("C5").Value = '2'
$ws.Rows.Item(5).AutoFit()

$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '0'
$ws.Range("C6").Value = '0'
$ws.Rows.Item(6).AutoFit()

$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '30 h'
$ws.Range("C7").Value = '30 h'
$ws.Rows.Item(7).AutoFit()

$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2022'
$ws.Range("C8").Value = '01/01/2022'
$ws.Rows.Item(8).AutoFit()

$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EF-6,EM-4'
$ws.Range("C9").Value = 'EF-6,EM-4'
$ws.Rows.Item(9).AutoFit()

$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Possibilitar ao estudante de Engenharia de Materiais o acesso a ferramentas computacionais modernas, de modo a que consiga descrever e quantificar conceitos vistos em outras disciplinas, como Ciência dos Materiais, Diagramas de Fases, Cinética de Transformação em Materiais, Termodinâmica, Propriedades Elétricas, Magnéticas, Térmicas e Ópticas, etc. Ao final do curso, o aluno será capaz de aplicar e entender resultados de simulações computacionais realistas aplicadas a diversas classes de materiais.'
$ws.Range("C10").Value = 'Possibilitar ao estudante de Engenharia de Materiais o acesso a ferramentas computacionais modernas, de modo a que consiga descrever e quantificar conceitos vistos em outras disciplinas, como Ciência dos Materiais, Diagramas de Fases, Cinética de Transformação em Materiais, Termodinâmica, Propriedades Elétricas, Magnéticas, Térmicas e Ópticas, etc. Ao final do curso, o aluno será capaz de aplicar e entender resultados de simulações computacionais realistas aplicadas a diversas classes de materiais.'
$ws.Rows.Item(10).RowHeight = 60

$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'Provide to Materials Engineering students access to modern computational tools, so that they can describe and quantify concepts seen in other disciplines, such as Materials Science, Phase Diagrams, Transformation Kinetics in Materials, Thermodynamics, Electrical, Magnetic, Thermal and Optical Properties, etc. At the end of the course, the student will be able to apply and understand the results of realistic computer simulations applied to different classes of materials.'
$ws.Range("C11").Value = 'Provide to Materials Engineering students access to modern computational tools, so that they can describe and quantify concepts seen in other disciplines, such as Materials Science, Phase Diagrams, Transformation Kinetics in Materials, Thermodynamics, Electrical, Magnetic, Thermal and Optical Properties, etc. At the end of the course, the student will be able to apply and understand the results of realistic computer simulations applied to different classes of materials.'
$ws.Rows.Item(11).RowHeight = 60

$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Rows.Item(12).AutoFit()

$ws.Range("B13").Value = '3480026 - João Paulo Pascon'
$ws.Range("C13").Value = '3480026 - João Paulo Pascon'
$ws.Rows.Item(13).AutoFit()

$ws.Range("B14").Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Range("C14").Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Rows.Item(14).AutoFit()

$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B15").Value = 'Tratamento de imagens em materialografia; Ajuste de equações empíricas ; Potenciais interatômicos e dinâmica molecular clássica; Descrição da Cinética de nucleação e crescimento; Método dos Elementos Finitos; Métodos de Monte Carlo; Crescimento de grão; Cálculo de Diagramas de fases.'
$ws.Range("C15").Value = 'Tratamento de imagens em materialografia; Ajuste de equações empíricas ; Potenciais interatômicos e dinâmica molecular clássica; Descrição da Cinética de nucleação e crescimento; Método dos Elementos Finitos; Métodos de Monte Carlo; Crescimento de grão; Cálculo de Diagramas de fases.'
$ws.Rows.Item(15).RowHeight = 60

$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("B16").Value = 'Image processing in materialography; Adjusting empirical equations; Interatomic potentials and classical molecular dynamics; Description of nucleation and growth kinetics; Finite Element Method; Monte Carlo methods; Grain growth; Calculation of phase diagrams.'
$ws.Range("C16").Value = 'Image processing in materialography; Adjusting empirical equations; Interatomic potentials and classical molecular dynamics; Description of nucleation and growth kinetics; Finite Element Method; Monte Carlo methods; Grain growth; Calculation of phase diagrams.'
$ws.Rows.Item(16).RowHeight = 60

$ws.Range("A17").Value = 'Programa:'
$ws.Range("B17").Value = '- Tratamento de imagens: resolução, definição, contraste, saturação; uso de técnicas automatizadas de determinação de tamanho e distribuição de partículas.- Proposição e ajuste de equações empíricas a resultados de medidas experimentais: as diversas propostas de relações para a deformação plástica e encruamento.- Potenciais interatômicos e o método de dinâmica molecular clássica; simulação de solidificação de um metal puro.- Cinética de nucleação e crescimento: a equação de Johnson-Mehl-Avrami-Kolmogorov (JMAK) e sua aplicação computacional.- Elementos finitos: estudo do estado de tensão de materiais sob carregamentos mecânicos; simulação de transferência de calor em tratamentos térmicos.- Método de Monte Carlo aplicado à transição ferro-paramagnética e à cinética de crescimento de grão- Cálculo de diagramas de fases: curvas de energia livre, o método CALPHAD; Thermo-Calc e Dictra.'
$ws.Range("C17").Value = '- Tratamento de imagens: resolução, definição, contraste, saturação; uso de técnicas automatizadas de determinação de tamanho e distribuição de partículas.- Proposição e ajuste de equações empíricas a resultados de medidas experimentais: as diversas propostas de relações para a deformação plástica e encruamento.- Potenciais interatômicos e o método de dinâmica molecular clássica; simulação de solidificação de um metal puro.- Cinética de nucleação e crescimento: a equação de Johnson-Mehl-Avrami-Kolmogorov (JMAK) e sua aplicação computacional.- Elementos finitos: estudo do estado de tensão de materiais sob carregamentos mecânicos; simulação de transferência de calor em tratamentos térmicos.- Método de Monte Carlo aplicado à transição ferro-paramagnética e à cinética de crescimento de grão- Cálculo de diagramas de fases: curvas de energia livre, o método CALPHAD; Thermo-Calc e Dictra.'
$ws.Rows.Item(17).RowHeight = 120

$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("B18").Value = '- Image treatment: resolution, definition, contrast, saturation; use of automated techniques for determining particle size and distribution.- Proposition and fit of empirical equations to results of experimental measures: the various proposals for relationships for plastic deformation and hardening.- Interatomic potentials and the classical molecular dynamics method; simulation of solidification of a pure metal.- Nucleation and growth kinetics: the Johnson-Mehl-Avrami-Kolmogorov (JMAK) equation and its computational application.- Finite element method: study of the stress state of materials under mechanical loads; simulation of heat transfer applied to heat treatments.- Monte Carlo method applied to the ferro-paramagnetic transition and to grain growth kinetics- Calculation of phase diagrams: free energy curves, the CALPHAD method; Thermo-Calc and Dictra.'
$ws.Range("C18").Value = '- Image treatment: resolution, definition, contrast, saturation; use of automated techniques for determining particle size and distribution.- Proposition and fit of empirical equations to results of experimental measures: the various proposals for relationships for plastic deformation and hardening.- Interatomic potentials and the classical molecular dynamics method; simulation of solidification of a pure metal.- Nucleation and growth kinetics: the Johnson-Mehl-Avrami-Kolmogorov (JMAK) equation and its computational application.- Finite element method: study of the stress state of materials under mechanical loads; simulation of heat transfer applied to heat treatments.- Monte Carlo method applied to the ferro-paramagnetic transition and to grain growth kinetics- Calculation of phase diagrams: free energy curves, the CALPHAD method; Thermo-Calc and Dictra.'
$ws.Rows.Item(18).RowHeight = 120

$ws.Range("A19").Value = 'Avaliação:'
$ws.Rows.Item(19).AutoFit()

$ws.Range("A20").Value = 'Método:'
$ws.Range("B20").Value = 'Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados. Trabalho baseado em Projeto'
$ws.Range("C20").Value = 'Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados. Trabalho baseado em Projeto'
$ws.Rows.Item(20).RowHeight = 60

$ws.Range("A21").Value = 'Critério:'
$ws.Range("B21").Value = 'Média aritmética de trabalhos propostos ao longo do curso (60%) e do Trabalho final em grupo (40%).'
$ws.Range("C21").Value = 'Média aritmética de trabalhos propostos ao longo do curso (60%) e do Trabalho final em grupo (40%).'
$ws.Rows.Item(21).RowHeight = 60

$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B22").Value = 'Não haverá exame de recuperação.'
$ws.Range("C22").Value = 'Não haverá exame de recuperação.'
$ws.Rows.Item(22).RowHeight = 60

$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("B23").Value = '- Richard LESAR, Computational Materials Science – Fundamentals to Applications. MRS, 2013.- Rob Phillips, Crystals, Defects and Microstructures – Modelling across scales. Cambridge, 2001.- Artigos publicados em revistas como Computational Materials Science, Calphad, Journal of Alloys and Compounds, etc.'
$ws.Range("C23").Value = '- Richard LESAR, Computational Materials Science – Fundamentals to Applications. MRS, 2013.- Rob Phillips, Crystals, Defects and Microstructures – Modelling across scales. Cambridge, 2001.- Artigos publicados em revistas como Computational Materials Science, Calphad, Journal of Alloys and Compounds, etc.'
$ws.Rows.Item(23).RowHeight = 120

$ws.Range("A24").Value = 'Requisitos:'
$ws.Rows.Item(24).AutoFit()

$ws.Range("B25").Value = 'LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)
'
$ws.Range("C25").Value = 'LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)
'
$ws.Rows.Item(25).RowHeight = 30

